# NIT-9016458754 Estado de Cuenta (EC) update
#
# The worker "RENE PORRAS PUELLO" previously had two overdue-period rows
# (2507 and 2506). This update removes the 2507 row entirely (its data was
# a duplicate/placeholder) and keeps a single row for period 2506, then
# refreshes the summary totals ("VALOR MORA" and "Cant. Periodos") to match
# the now-single remaining period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "2507" detail row (row 17) entirely - this shifts every
# row below it up by one (so the old rows 22/23 signature block becomes
# rows 21/22).
$ws.Rows("17").Delete()

# The remaining detail row (now the only one, row 16) represents period
# 2506.
$ws.Range("E16").Value = "2506"

# Refresh the summary header: total "VALOR MORA" and "Cant. Periodos" now
# reflect the single remaining period/row.
$ws.Range("E11").Value = 56940
$ws.Range("F13").Value = 1
